$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E11 gets a new value (master_worker_response_tracke_ruan_r1_Sep-30-2023.csv)
$ws.Range("E11").Value = "master_worker_response_tracke_ruan_r1_Sep-30-2023.csv"

# Insert one new row before old row 12, shifting old rows 12-13 down to 13-14
$ws.Rows("12:12").Insert()

# New row 12: RD's Sep-30 submission
$ws.Range("A12").Value = "pairwise"
$ws.Range("B12").Value = "Sept-30-2023"
$ws.Range("C12").Value = "RD"
$ws.Range("D12").Value = "all_submitted_tracker_ruan_r2_Sep-30-2023.csv"

# New row 15 (at the end): SB's Oct-1 submission
$ws.Range("A15").Value = "paiewise"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Oct-1-2023"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "SB"
$ws.Range("D15").Value = "all_submitted_tracker_SB_Oct-01-2023.csv"

$ws.Range("D15").Select()
